$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.03"
$ws.Range("E2").Value = "'-0.83%"

$ws.Range("D3").Value = "'35.80"
$ws.Range("E3").Value = "'-0.21%"

$ws.Range("D4").Value = "'5.048"
$ws.Range("E4").Value = "'-0.34%"

$ws.Range("D5").Value = "'0.07978"
$ws.Range("E5").Value = "'-1.21%"

$ws.Range("D6").Value = "'1.852"
$ws.Range("E6").Value = "'-5.14%"

$ws.Range("D7").Value = "'7.776"
$ws.Range("E7").Value = "'-0.04%"

$ws.Range("D8").Value = "'0.9222"
$ws.Range("E8").Value = "'-0.62%"

$ws.Range("D9").Value = "'0.1292"
$ws.Range("E9").Value = "'-5.14%"

$ws.Range("D10").Value = "'0.1885"
$ws.Range("E10").Value = "'-0.89%"

$ws.Range("D11").Value = "'0.09091"
$ws.Range("E11").Value = "'-1.72%"

$ws.Range("D12").Value = "'0.03428"
$ws.Range("E12").Value = "'-2.91%"

$ws.Range("D13").Value = "'0.09867"
$ws.Range("E13").Value = "'-0.07%"

$ws.Range("D14").Value = "'0.001405"
$ws.Range("E14").Value = "'-1.85%"

$ws.Range("D15").Value = "'0.006190"
$ws.Range("E15").Value = "'6.45%"

$ws.Range("D16").Value = "'3.858"
$ws.Range("E16").Value = "'7.99%"

$ws.Range("D17").Value = "'4.123"
$ws.Range("E17").Value = "'-0.83%"

$ws.Range("E18").Value = "'14.06%"

$ws.Range("D19").Value = "'0.3404"
$ws.Range("E19").Value = "'-1.23%"

$ws.Range("D20").Value = "'0.1310"
$ws.Range("E20").Value = "'-0.98%"

$ws.Range("D21").Value = "'4.811"
$ws.Range("E21").Value = "'-1.72%"

$ws.Range("D22").Value = "'0.2499"

$ws.Range("D23").Value = "'0.04427"
$ws.Range("E23").Value = "'0.77%"

$ws.Range("D24").Value = "'0.001235"
$ws.Range("E24").Value = "'1.12%"

$ws.Range("D25").Value = "'0.004871"
$ws.Range("E25").Value = "'2.10%"

$ws.Range("D27").Value = "'0.0001304"
$ws.Range("E27").Value = "'-21.01%"

$ws.Range("E28").Value = "'42.08%"

$ws.Range("E39").Value = "'-1.42%"

$ws.Range("D40").Value = "'0.05150"
$ws.Range("E40").Value = "'2.84%"

$ws.Range("D41").Value = "'0.007536"
$ws.Range("E41").Value = "'-0.80%"

$ws.Range("D42").Value = "'0.01013"
$ws.Range("E42").Value = "'-8.74%"

$ws.Range("D43").Value = "'0.1351"
$ws.Range("E43").Value = "'-1.96%"

$ws.Range("D44").Value = "'0.002116"
$ws.Range("E44").Value = "'0.92%"

$ws.Range("D45").Value = "'0.009887"
$ws.Range("E45").Value = "'-8.50%"

$ws.Range("D46").Value = "'0.00006170"
$ws.Range("E46").Value = "'-3.80%"

$ws.Range("E47").Value = "'0.00%"

$ws.Range("D48").Value = "'63.66"
$ws.Range("E48").Value = "'-2.00%"

$ws.Range("D49").Value = "'0.001251"
$ws.Range("E49").Value = "'5.01%"

$ws.Range("E50").Value = "'0.00%"

$ws.Range("E51").Value = "'0.00%"
